# edit.ps1 - applies the tracked changes described in the commit diff
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $find)
    }
    return $ok
}

# 1. Title / course name (paragraph near top of doc)
Replace-Text "CURSO PREPARATÓRIO Prep. Básico PBH-SMS- 2025 COM GRAVAÇÃO DAS AULAS." "CURSO PREPARATÓRIO PREPARATÓRIO CONTEÚDO BÁSICO NÍVEL SUPERIOR - PBH/SMS COM GRAVAÇÃO DAS AULAS." | Out-Null

# 2. Student name everywhere it appears with a trailing space (Nome: field + both signature blocks)
Replace-Text "Leonardo de Oliveira " "Leonardo de Oliveira Narciso" | Out-Null

# 3. CEP
Replace-Text "CEP: 34011060" "CEP: 30180060" | Out-Null

# 4. Cidade/UF
Replace-Text "Cidade/UF: - MG" "Cidade/UF: Belo Horizonte" | Out-Null

# 5. Forma de Pagamento
Replace-Text "Forma de Pagamento: crédito " "Forma de Pagamento: Cartão de crédito - 10x " | Out-Null

# 6. Total de cada parcela
Replace-Text "Total de cada parcela: 345,00" "Total de cada parcela: 100" | Out-Null

# 7. Total do contrato
Replace-Text "Total do contrato:  R$ 1250,00" "Total do contrato:  R$ 1000" | Out-Null

# 8. Clausula 1a - objeto do contrato
Replace-Text "Cláusula 1ª. Constitui objeto do presente Contrato Prep. Básico PBH-SMS- 2025 COM GRAVAÇÃO DAS AULAS." "Cláusula 1ª. Constitui objeto do presente Contrato PREPARATÓRIO CONTEÚDO BÁSICO NÍVEL SUPERIOR - PBH/SMS COM GRAVAÇÃO DAS AULAS." | Out-Null

# 9. Presencial / carga horaria
Replace-Text " Presencial com acesso a aulas gravadas em sala a ser ofertado pela CONTRATADA em favor do CONTRATANTE, com carga horária total oferecida de APROXIMADAMENTE 60 horas em periodos de 30 min." " Presencial com acesso a aulas gravadas em sala a ser ofertado pela CONTRATADA em favor do CONTRATANTE, com carga horária total oferecida de APROXIMADAMENTE 60." | Out-Null

# 10. Paragrafo Primeiro - datas e dias da semana
Replace-Text "§ Primeiro: As aulas serão ministradas presencialmente, à Rua Juiz de Fora 231 Barro Preto no período de 09 de abril de 2025 com término previsto para dia 30 de agosto de 2025, seg à sexta  de 19:00h as 22:15h e aos seg à sexta de 09:00h às 12:15h. " "§ Primeiro: As aulas serão ministradas presencialmente, à Rua Juiz de Fora 231 Barro Preto no período de 05 de maio de 2025 com término previsto para dia 20 de agosto de 2025, Segundas e Quintas  de 19:00h as 22:15h e aos Segundas e Quintas de 09:00h às 12:15h. " | Out-Null

# 11. Split the run containing "...falha técnica.§ Terceiro: " into three runs:
#     "...falha técnica." | " " | "§Terceiro: "
$r11 = $d.Content
$ok11 = $r11.Find.Execute("falha técnica.§ Terceiro: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok11) {
    $start11 = $r11.Start
    $r11.Text = "falha técnica."
    $afterFirst = $start11 + ("falha técnica.").Length
    $p11a = $d.Range($afterFirst, $afterFirst)
    $p11a.InsertAfter(" ")
    $afterSpace = $afterFirst + 1
    $p11b = $d.Range($afterSpace, $afterSpace)
    $p11b.InsertAfter("§Terceiro: ")
    # force the three inserted spans into distinct runs (no net formatting change)
    $spaceRange = $d.Range($afterFirst, $afterSpace)
    $spaceRange.Bold = 1
    $spaceRange.Bold = 0
    $thirdRange = $d.Range($afterSpace, $afterSpace + ("§Terceiro: ").Length)
    $thirdRange.Bold = 1
    $thirdRange.Bold = 0
} else {
    Write-Output "NOT FOUND: falha técnica.§ Terceiro: "
}

# 12. Split off a new empty paragraph before the "Belo Horizonte..." signature-date paragraph,
#     and shorten/replace its text.
$r12 = $d.Content
$ok12 = $r12.Find.Execute("Belo Horizonte, Belo Horizonte, 09 de abril de 2025 de {{mes_contrato}} de {{ano_contrato}}.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok12) {
    $r12.InsertParagraphBefore()
    $ok12b = $d.Content.Find.Execute("Belo Horizonte, Belo Horizonte, 09 de abril de 2025 de {{mes_contrato}} de {{ano_contrato}}.", $true, $false, $false, $false, $false, $true, 1, $false, "Belo Horizonte, 11 de abril de 2025", 2)
    if (-not $ok12b) {
        Write-Output "NOT FOUND (pass 2): Belo Horizonte, Belo Horizonte..."
    }
} else {
    Write-Output "NOT FOUND: Belo Horizonte, Belo Horizonte..."
}

Write-Output "done"
